# Planlegging.docx update:
#  - Re-order / extend the "Sider:" paragraph: move "Om Oss, " in front of "Video",
#    and replace the "Spill nettside" list item with a short summary sentence
#    ("Video og Spill. 4 nettsider totalt.").
#  - Add a new paragraph right after it describing where the sketches were found.

$d = $word.ActiveDocument

# 1) Rewrite the tail of the "Sider:" paragraph.
#    Before: "Sider: Hjemme side, Video, Spill nettside, Om Oss, "
#    After : "Sider: Hjemme side, Om Oss, Video og Spill. 4 nettsider totalt. "
$replaced = $d.Content.Find.Execute(
    "Video, Spill nettside, Om Oss, ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Om Oss, Video og Spill. 4 nettsider totalt. ", 2)

if (-not $replaced) {
    Write-Host "WARNING: could not find the 'Sider:' list tail to replace"
}

# 2) Locate that paragraph (by its now-updated text) and add a new paragraph
#    right after it with the "Skisser funnet i ..." sentence.
$siderIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Sider: Hjemme side, Om Oss, Video og Spill.")) {
        $siderIndex = $i
        break
    }
}

if ($siderIndex -gt 0) {
    $siderPara = $d.Paragraphs.Item($siderIndex)
    $siderPara.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($siderIndex + 1)
    $newRange = $newPara.Range
    $newRange.Collapse(1)
    $newRange.InsertBefore("Skisser funnet i skisse mappen. ")
} else {
    Write-Host "WARNING: could not locate the updated 'Sider:' paragraph"
}
